$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "test"
$ws.Name = "test"

# Re-apply the "Normal" style to the used range (A1:D2) so that the
# cells pick up an explicit style definition (font/alignment/protection
# applied) instead of relying on the implicit default style.
$rng = $ws.Range("A1:D2")
$rng.Style = "Normal"
